$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "ödeme"
$ws.Range("B44").Value = "Müşterinin ödediği tutarlları görmek istiyorum."
$ws.Range("D44").Value = "Salesforce da plakayı yazıp ara. Ödeme - Evrak - Teslimat- Kazanıldı aşamasında olan kaydı seç. Ödeme sekmesini seç. En altta müşterinin yapmış olduğu ödemeleri görebilirsin."
$ws.Range("C44").Value = "Müşterinin ödediği tutarları ödeme sekmesinden görebilirsin."
$ws.Range("E44").Value = "Product Manager"
$ws.Range("F44").Value = "ödeme.JPG"

$ws.Range("F45").Select()
